# Fix header labels on the existing sheets (Weekly/Monthly quantity columns
# renamed to match the forecasting pipeline's expected column names).
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Weekly Quantity")
$ws2 = $wb.Worksheets.Item("Monthly Trend")

$ws1.Range("B1").Value = "Weekly_PO_Qty"
$ws2.Range("B1").Value = "Monthly_PO_Qty"

# Add the new "PO Forecast" sheet after "Monthly Trend" so it becomes the
# third/last tab.
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "PO Forecast"

# Header row
$ws3.Range("A1").Value = "ds"
$ws3.Range("B1").Value = "PO_Forecast"
$ws3.Range("C1").Value = "yhat_lower"
$ws3.Range("D1").Value = "yhat_upper"

# Forecast data rows (ds, PO_Forecast, yhat_lower, yhat_upper)
$csv = @"
44934.99999999999,228,-410.4932435600022,859.1393161235943
44941.99999999999,233,-376.8947906502653,847.2370996090932
44948.99999999999,239,-399.6392454518359,860.0218364442676
44955.99999999999,244,-381.7531993866939,821.7563216302245
44962.99999999999,250,-403.398519983331,879.6843019837203
44969.99999999999,255,-413.1287705968193,879.7002737064048
44976.99999999999,261,-340.7637337073002,915.59855428891
44983.99999999999,266,-351.4146721376774,846.9609611665835
45004.99999999999,283,-357.8984721033828,852.4895613365679
45011.99999999999,288,-389.5696100580965,907.0781705592138
45025.99999999999,299,-300.5922159342132,916.1149699857464
45032.99999999999,305,-299.1106973967053,957.5740644839062
45039.99999999999,310,-291.076998370101,934.7047624454884
45046.99999999999,315,-322.9771889756356,926.9582226858903
45060.99999999999,326,-303.0677014798375,993.004833857039
45067.99999999999,332,-309.2029652543403,907.5999198525179
45074.99999999999,337,-318.9786037022743,965.0677520674819
45095.99999999999,354,-241.3537652058225,949.9127686344829
45102.99999999999,359,-276.5544165592271,981.0856265905377
45109.99999999999,365,-218.8753503163042,1050.358598517618
45116.99999999999,370,-246.864135510122,1002.069906323047
45123.99999999999,376,-271.2562390508762,928.6940351272958
45137.99999999999,387,-231.3989733407674,1013.889571678055
45165.99999999999,409,-214.0943130904239,1025.751494920505
45193.99999999999,431,-226.7688532106395,1061.305167272373
45214.99999999999,447,-186.4886603240689,1024.262888725233
45221.99999999999,453,-182.3412301880074,1077.782603227748
45228.99999999999,458,-171.1825864193552,1051.836447567895
45235.99999999999,464,-172.4165668869585,1120.739896359273
45242.99999999999,469,-138.6136284959477,1092.356418513674
45256.99999999999,480,-152.9237819019959,1118.455652343468
45270.99999999999,491,-142.2042531151615,1111.67423530389
45277.99999999999,497,-87.25163535151955,1140.42729225696
45298.99999999999,513,-97.84159631120296,1128.712031602635
45305.99999999999,519,-106.5443012715124,1199.048500865504
45312.99999999999,524,-40.11841710297784,1130.997271100374
45319.99999999999,530,-83.40498199862375,1156.829526160674
45326.99999999999,535,-91.60573605448347,1144.45395173004
45333.99999999999,541,-115.0290018007539,1172.102298992351
45340.99999999999,546,-93.02208397554803,1207.508233243146
45347.99999999999,552,-55.19619359096368,1159.046920515786
45354.99999999999,557,-81.05790545883364,1188.502959506497
45361.99999999999,563,-60.72213639026643,1188.604523024591
45375.99999999999,574,-7.800208581160937,1181.668161218981
45389.99999999999,585,-54.33364589716685,1178.061481133879
45396.99999999999,590,-51.7354168869799,1206.116253578567
45403.99999999999,596,-72.06602418745426,1196.281446956017
45410.99999999999,601,-2.434364600724424,1233.334050411602
45417.99999999999,607,-55.8458141315491,1213.934511114643
45466.99999999999,645,0.9821262356451542,1284.806746592654
45473.99999999999,651,23.18769291766399,1263.358875121585
45480.99999999999,656,23.38932140678098,1274.694879022434
45487.99999999999,662,7.581338155361554,1278.158801975783
45529.99999999999,695,93.1399316590619,1297.021447138035
45536.99999999999,700,66.37407351135698,1330.009347472816
45543.99999999999,706,56.40647976185933,1371.540057312838
45557.99999999999,717,92.28074113428376,1330.779108312222
45564.99999999999,722,128.0798912592305,1341.064507002274
45585.99999999999,739,155.126520367054,1337.417918431052
45599.99999999999,750,89.85116536111737,1380.8381037471
45606.99999999999,755,136.5126730233247,1411.320125484511
45613.99999999999,761,159.4287331477285,1357.823969289018
45641.99999999999,783,132.72712680059,1399.50261514568
45648.99999999999,788,120.6724736361031,1427.635371625442
45655.99999999999,794,160.5687780261879,1443.200924840191
45662.99999999999,799,195.7355193533761,1418.162751249637
45669.99999999999,805,206.9822573022963,1450.623310036378
45676.99999999999,810,219.6235194680325,1433.938290606466
45683.99999999999,816,220.9827374933692,1453.227707019622
45690.99999999999,821,182.6756279685951,1404.488828983327
45697.99999999999,827,196.5371999939747,1396.652930856747
"@

$lines = $csv -split "`r?`n"
$r = 2
foreach ($line in $lines) {
    if ($line.Trim().Length -eq 0) { continue }
    $vals = $line -split ","
    $ws3.Cells.Item($r, 1).Value = [double]$vals[0]
    $ws3.Cells.Item($r, 2).Value = [double]$vals[1]
    $ws3.Cells.Item($r, 3).Value = [double]$vals[2]
    $ws3.Cells.Item($r, 4).Value = [double]$vals[3]
    $r = $r + 1
}

# Match formatting to the rest of the workbook: header style (bold, centered,
# bordered) and date-serial number format for column A, reusing the exact
# same style entries already present in styles.xml (avoids duplicating them).
$ws1.Range("A1:B1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)

$ws1.Range("A2").Copy()
$ws3.Range("A2:A72").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Restore the original active sheet/selection (adding the sheet switches to
# it by default).
$ws1.Activate()
[void]$ws1.Range("A1").Select()
